$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 00:08"

# --- Reorder Cuba / Suazilandia: Suazilandia now precedes Cuba ---
$ws.Cells.Item(119, 1).Value = "Suazilandia"
$ws.Cells.Item(120, 1).Value = "Cuba"

# --- Refresh per-country daily statistics: Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 4697704
$ws.Cells.Item(4, 3).Value = 62719
$ws.Cells.Item(4, 4).Value = 2318773
$ws.Cells.Item(4, 5).Value = 2222621
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 1025
$ws.Cells.Item(4, 8).Value = 156310

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 2662485
$ws.Cells.Item(5, 3).Value = 48696
$ws.Cells.Item(5, 4).Value = 1844051
$ws.Cells.Item(5, 5).Value = 725959
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 1098
$ws.Cells.Item(5, 8).Value = 92475

# Row 15: Colombia
$ws.Cells.Item(15, 2).Value = 295508
$ws.Cells.Item(15, 3).Value = 9488
$ws.Cells.Item(15, 4).Value = 154387
$ws.Cells.Item(15, 5).Value = 131016
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 295
$ws.Cells.Item(15, 8).Value = 10105

# Row 21: Alemania
$ws.Cells.Item(21, 2).Value = 210665
$ws.Cells.Item(21, 3).Value = 1012
$ws.Cells.Item(21, 4).Value = 192300
$ws.Cells.Item(21, 5).Value = 9141
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 3
$ws.Cells.Item(21, 8).Value = 9224

# Row 28: Egipto
$ws.Cells.Item(28, 2).Value = 94078
$ws.Cells.Item(28, 3).Value = 321
$ws.Cells.Item(28, 4).Value = 39638
$ws.Cells.Item(28, 5).Value = 49635
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 31
$ws.Cells.Item(28, 8).Value = 4805

# Row 52: Barein
$ws.Cells.Item(52, 2).Value = 40982
$ws.Cells.Item(52, 3).Value = 227
$ws.Cells.Item(52, 4).Value = 37840
$ws.Cells.Item(52, 5).Value = 2995
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 147

# Row 55: Kirguistan
$ws.Cells.Item(55, 2).Value = 35805
$ws.Cells.Item(55, 3).Value = 582
$ws.Cells.Item(55, 4).Value = 25037
$ws.Cells.Item(55, 5).Value = 9390
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 14
$ws.Cells.Item(55, 8).Value = 1378

# Row 76: Costa de Marfil
$ws.Cells.Item(76, 2).Value = 16047
$ws.Cells.Item(76, 3).Value = 69
$ws.Cells.Item(76, 4).Value = 11428
$ws.Cells.Item(76, 5).Value = 4517
$ws.Cells.Item(76, 6).Value = 0
$ws.Cells.Item(76, 7).Value = 2
$ws.Cells.Item(76, 8).Value = 102

# Row 96: Mauritania
$ws.Cells.Item(96, 2).Value = 6310
$ws.Cells.Item(96, 3).Value = 15
$ws.Cells.Item(96, 4).Value = 4962
$ws.Cells.Item(96, 5).Value = 1191
$ws.Cells.Item(96, 6).Value = 0
$ws.Cells.Item(96, 7).Value = 0
$ws.Cells.Item(96, 8).Value = 157

# Row 101: Republica de Yibuti
$ws.Cells.Item(101, 2).Value = 5084
$ws.Cells.Item(101, 3).Value = 3
$ws.Cells.Item(101, 4).Value = 4999
$ws.Cells.Item(101, 5).Value = 27
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 58

# Row 102: Republica de Africa Central
$ws.Cells.Item(102, 2).Value = 4608
$ws.Cells.Item(102, 3).Value = 3
$ws.Cells.Item(102, 4).Value = 1606
$ws.Cells.Item(102, 5).Value = 2943
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 59

# Row 119: Suazilandia
$ws.Cells.Item(119, 2).Value = 2648
$ws.Cells.Item(119, 3).Value = 71
$ws.Cells.Item(119, 4).Value = 1214
$ws.Cells.Item(119, 5).Value = 1393
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 41

# Row 120: Cuba
$ws.Cells.Item(120, 2).Value = 2608
$ws.Cells.Item(120, 3).Value = 11
$ws.Cells.Item(120, 4).Value = 2355
$ws.Cells.Item(120, 5).Value = 166
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 7).Value = 0
$ws.Cells.Item(120, 8).Value = 87

# Row 129: Ruanda
$ws.Cells.Item(129, 2).Value = 2022
$ws.Cells.Item(129, 3).Value = 28
$ws.Cells.Item(129, 4).Value = 1106
$ws.Cells.Item(129, 5).Value = 911
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 5

# Row 138: Tunez
$ws.Cells.Item(138, 2).Value = 1535
$ws.Cells.Item(138, 3).Value = 21
$ws.Cells.Item(138, 4).Value = 1195
$ws.Cells.Item(138, 5).Value = 290
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 0
$ws.Cells.Item(138, 8).Value = 50

# Row 149: Republica del Chad
$ws.Cells.Item(149, 2).Value = 936
$ws.Cells.Item(149, 3).Value = 1
$ws.Cells.Item(149, 4).Value = 813
$ws.Cells.Item(149, 5).Value = 48
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 75

# Row 209: San Cristobal y Nieves
$ws.Cells.Item(209, 2).Value = 17
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 16
$ws.Cells.Item(209, 5).Value = 1
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0
